{"js": "// Colors the two \"to-do\" lines about the win screen / hiding enemy pieces\n// in green (92D050), matching the shading already used for the other\n// to-do items (\"Validera att alla enheter...\" / \"Fixa r\u00f6relselogiken\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst GREEN = \"#92D050\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text.indexOf(\"vinst sk\\u00e4rm\") !== -1 || text.indexOf(\"fieendepj\\u00e4serna\") !== -1) {\n    para.font.color = GREEN;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Colors the two \"to-do\" lines about the win screen / hiding enemy pieces\n# in green (RGB 92D050 -> BGR long 5296274), matching the shading already\n# used for the other to-do items (\"Validera att alla enheter...\" /\n# \"Fixa r\u00f6relselogiken\").\n$d = $word.ActiveDocument\n$green = 5296274  # RGB(0x92,0xD0,0x50) packed as BGR for wdColor\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*vinst sk*\" -or $t -like \"*fieendepj*\") {\n        $p.Range.Font.Color = $green\n    }\n}\n"}
